$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous layout (A1 had "Topic Name", B3 had "Name")
$ws.Range("A1").ClearContents()

# Row 2: header row for the first (Name) table
$ws.Range("A2").Value = "Topic Name"
$ws.Range("B2").Value = "Quick Notes"

# Row 3: "Name" column header
$ws.Range("B3").Value = "Name"

# Rows 4-7: names under the first topic
$ws.Range("B4").Value = "Avinash Jadhav"
$ws.Range("B5").Value = "Sanket Chor"
$ws.Range("B6").Value = "Harshwardhan Pachoute"
$ws.Range("B7").Value = "Om Dighe"

# Row 9: header row for the second (Musical Instrument) table
$ws.Range("A9").Value = "Topic Name"
$ws.Range("B9").Value = "Musical Instrument "

# Rows 10-13: remaining names
$ws.Range("B10").Value = "Pranit Vichare"
$ws.Range("B11").Value = "Priyanka Rasal"
$ws.Range("B12").Value = "Tajas Sutar"
$ws.Range("B13").Value = "Deepak Dixit"

# Column B width to fit the longer strings (closest achievable to 28.453125)
$ws.Columns.Item(2).ColumnWidth = 27.6

# Zoom / selection state from the saved view
$ws.Application.ActiveWindow.Zoom = 190
$ws.Range("C9").Select()
